$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.794.58'
$ws.Range("E2").Value = '  -0.63%  '
$ws.Range("D3").Value = '3.000.50'
$ws.Range("E3").Value = '  -0.04%  '
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '598.32'
$ws.Range("E5").Value = '  +2.64%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.31'
$ws.Range("E6").Value = '  +0.29%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.998'
$ws.Range("E7").Value = '  -0.18%  '
$ws.Range("D8").Value = '3.000.28'
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.514'
$ws.Range("E9").Value = '  -2.02%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.149'
$ws.Range("E10").Value = '  +0.53%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.06'
$ws.Range("E11").Value = '  +4.54%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.455'
$ws.Range("E12").Value = '  -0.60%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000229'
$ws.Range("E13").Value = '  +0.19%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.34'
$ws.Range("E14").Value = '  -0.13%  '
$ws.Range("E15").Value = '  +3.49%  '
$ws.Range("D16").Value = '3.490.14'
$ws.Range("E16").Value = '  -0.26%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.94'
$ws.Range("E17").Value = '  -2.04%  '
$ws.Range("D18").Value = '61.739.23'
$ws.Range("E18").Value = '  -0.72%  '
$ws.Range("D19").Value = '2.999.83'
$ws.Range("E19").Value = '  -0.28%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '448.35'
$ws.Range("E20").Value = '  -2.59%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.01'
$ws.Range("E21").Value = '  +0.38%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.685'
$ws.Range("E22").Value = '  +0.08%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.35'
$ws.Range("E23").Value = '  -0.88%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '81.37'
$ws.Range("E24").Value = '  +0.08%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '10.96'
$ws.Range("E25").Value = '  +8.24%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.21'
$ws.Range("E26").Value = '  -0.52%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.09'
$ws.Range("E27").Value = '  -1.21%  '
$ws.Range("E28").Value = '  +0.10%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.71'
$ws.Range("E29").Value = '  +3.46%  '
$ws.Range("B30").Value = 'NEARProtocol'
$ws.Range("C30").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.31'
$ws.Range("E30").Value = '  +3.65%  '
$ws.Range("B31").Value = 'FirstDigitalUSD'
$ws.Range("C31").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.00'
$ws.Range("E31").Value = '  -0.12%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.07'
$ws.Range("E32").Value = '  -1.47%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '27.25'
$ws.Range("E33").Value = '  -3.56%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.110'
$ws.Range("E34").Value = '  +1.90%  '
$ws.Range("D35").Value = '0.0₃0830'
$ws.Range("E35").Value = '  +4.27%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.02'
$ws.Range("E36").Value = '  -0.07%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.79'
$ws.Range("E37").Value = '  +0.76%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '50.43'
$ws.Range("E38").Value = '  +0.26%  '
$ws.Range("B39").Value = 'Stacks'
$ws.Range("C39").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.04'
$ws.Range("E39").Value = '  -3.15%  '
$ws.Range("B40").Value = 'Cosmos'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '8.96'
$ws.Range("E40").Value = '  -2.18%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.123'
$ws.Range("E41").Value = '  +7.93%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.88'
$ws.Range("E42").Value = '  +0.22%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '396.37'
$ws.Range("E43").Value = '  +0.74%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '40.45'
$ws.Range("E44").Value = '  +10.13%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.271'
$ws.Range("E45").Value = '  -0.57%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0352'
$ws.Range("E46").Value = '  -1.35%  '
$ws.Range("D47").Value = '2.697.35'
$ws.Range("E47").Value = '  -1.24%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '132.18'
$ws.Range("E48").Value = '  +2.36%  '
$ws.Range("E49").Value = '  +0.10%  '
$ws.Range("B50").Value = 'ThetaToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.17'
$ws.Range("E50").Value = '  -0.68%  '
$ws.Range("B51").Value = 'Stellar'
$ws.Range("C51").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.107'
$ws.Range("E51").Value = '  -1.57%  '

Write-Output "Updated cryptos list"
